# TradingModel_v2 - 2021/11/18 data updated
# Add a new daily TotalCapital data point (2021-11-18) to the bottom of the
# table on the active sheet, and move the "last row" date formatting from
# the previous last row (2021-11-17) to the newly appended row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The former last row (row 11, 2021-11-17) should now use the regular
# date/time number format used by all the other non-final rows.
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new data row (row 12) for 2021-11-18.
$ws.Range("A12").Value = 44518
$ws.Range("B12").Value = 65081.4

# The new last row gets the distinct "last row" date-only number format.
$ws.Range("A12").NumberFormat = "YYYY-MM-DD"
